# Updates odds/score-count figures on the "Jogos da Semana" sheet to match
# the latest FlashScore refresh (commit: "Atualizando o arquivo XLSX").
# Only specific cells in rows 9, 10, 14, 16, 17 and 22 change value; everything
# else in the workbook stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 ---
$ws.Range("G9").Value  = 1.67
$ws.Range("I9").Value  = 5.5
$ws.Range("AA9").Value = 19
$ws.Range("AE9").Value = 1.07
$ws.Range("AF9").Value = 9

# --- Row 10 ---
$ws.Range("AE10").Value = 1.02
$ws.Range("AG10").Value = 1.13
$ws.Range("AI10").Value = 1.69

# --- Row 14 ---
$ws.Range("AJ14").Value = 1.63

# --- Row 16 ---
$ws.Range("G16").Value = 2.22
$ws.Range("H16").Value = 3.25
$ws.Range("I16").Value = 2.82
$ws.Range("J16").Value = 1.85
$ws.Range("K16").Value = 1.75
$ws.Range("M16").Value = 2.5
$ws.Range("N16").Value = 6.7
$ws.Range("T16").Value = 9.5
$ws.Range("U16").Value = 5.6
$ws.Range("Y16").Value = 7.7

# --- Row 17 ---
$ws.Range("G17").Value  = 2.75
$ws.Range("H17").Value  = 3.3
$ws.Range("I17").Value  = 2.63
$ws.Range("J17").Value  = 2.1
$ws.Range("K17").Value  = 1.7
$ws.Range("P17").Value  = 10
$ws.Range("T17").Value  = 9
$ws.Range("Y17").Value  = 8
$ws.Range("AE17").Value = 1.06
$ws.Range("AF17").Value = 10
$ws.Range("AG17").Value = 1.33
$ws.Range("AH17").Value = 3.25

# --- Row 22 ---
$ws.Range("G22").Value  = 7
$ws.Range("H22").Value  = 3.8
$ws.Range("I22").Value  = 1.44
$ws.Range("J22").Value  = 1.87
$ws.Range("K22").Value  = 1.87
$ws.Range("N22").Value  = 17
$ws.Range("O22").Value  = 34
$ws.Range("P22").Value  = 21
$ws.Range("Q22").Value  = 81
$ws.Range("T22").Value  = 10
$ws.Range("V22").Value  = 19
$ws.Range("W22").Value  = 67
$ws.Range("Z22").Value  = 6.5
$ws.Range("AB22").Value = 10
$ws.Range("AE22").Value = 1.06
$ws.Range("AF22").Value = 10
$ws.Range("AI22").Value = 2
$ws.Range("AJ22").Value = 1.73
